$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first file (rows 2-3 share the same value)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 18:17:58"
$wsOverview.Range("G3").Value = "2016-08-29 18:17:58"

# "zh-cn" sheet: Priority column (ht -> mt), Correspond Handoff/Handback Datetime updates
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-29 18:17:53"
$wsZhCn.Range("H3").Value = "2016-08-29 18:17:53"
$wsZhCn.Range("K2").Value = "2016-08-29 18:18:22"
$wsZhCn.Range("K3").Value = "2016-08-29 18:18:22"

# "de-de" sheet: Correspond Handback Datetime update
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-29 18:18:29"
$wsDeDe.Range("K3").Value = "2016-08-29 18:18:29"
